$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (15) had a trailing empty-marker cell in column H
# (an inline string cell with no text) just past the real data in A:G.
# That trailing marker needs to move to the new last row (22) once the
# additional log entries below are appended.
$ws.Range("H15").ClearContents()

$newRows = @(
    @("5/30/2022 21:51", "Monday",  "Alcon", "Processed", "Sent",                "Alcon_Status Report_05302022.xlsx", "Alcon_Document Expiration Report_05302022.xlsx"),
    @("5/31/2022 09:46", "Tuesday", "Alcon", "Processed", "Sent",                "Alcon_Status Report_05312022.xlsx", "Alcon_Document Expiration Report_05312022.xlsx"),
    @("5/31/2022 11:48", "Tuesday", "Alcon", "Processed", "Send Mail Disabled",  "Alcon_Status Report_05312022.xlsx", "Alcon_Document Expiration Report_05312022.xlsx"),
    @("5/31/2022 11:53", "Tuesday", "Alcon", "Processed", "Send Mail Disabled",  "Alcon_Status Report_05312022.xlsx", "Alcon_Document Expiration Report_05312022.xlsx"),
    @("5/31/2022 12:01", "Tuesday", "Alcon", "Processed", "Send Mail Disabled",  "Alcon_Status Report_05312022.xlsx", "Alcon_Document Expiration Report_05312022.xlsx"),
    @("5/31/2022 12:33", "Tuesday", "Alcon", "Processed", "Send Mail Disabled",  "Alcon_Status Report_05312022.xlsx", "Alcon_Document Expiration Report_05312022.xlsx"),
    @("5/31/2022 12:48", "Tuesday", "Alcon", "Processed", "Send Mail Disabled",  "Alcon_Status Report_05312022.xlsx", "Alcon_Document Expiration Report_05312022.xlsx")
)

$r = 16
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# Restore the trailing empty-marker cell on the new last row.
$ws.Range("H22").Value = ""
